$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3, 8 and 9 are now empty/blank rows (data removed from source),
# so clear out their contents. The lattice-parameter cells (D:G) keep
# their existing numeric style but have no value.
foreach ($r in 3, 8, 9) {
    $ws.Range("A$r").Value = $null
    $ws.Range("B$r").Value = $null
    $ws.Range("C$r").Value = $null
    $ws.Range("D$r").Value = $null
    $ws.Range("E$r").Value = $null
    $ws.Range("F$r").Value = $null
    $ws.Range("G$r").Value = $null
    $ws.Range("H$r").Value = $null
    $ws.Range("I$r").Value = $null
}

# Update the active selection to A10, matching the new state after the
# blank rows were logged/skipped.
$ws.Range("A10").Select()
